# Update the division problems in the single 20x5 practice table.
# Several cells share identical "old" text (e.g. "89÷5=", "58÷4="), so a
# document-wide Find/Replace-All would clobber duplicates with the same
# value. Instead, target each table cell individually, scope the
# Find/Execute to that cell's Range, and use wdReplaceOne (1) rather than
# wdReplaceAll (2) -- this engine's Find.Execute replaces every match in
# the document regardless of range when ReplaceAll is requested, so
# ReplaceOne is required to keep the edit confined to a single cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# (row, column, old, new)
$edits = @(
    @(1, 1, "81÷8=", "34÷7="),
    @(1, 2, "83÷2=", "32÷9="),
    @(1, 3, "46÷4=", "96÷9="),
    @(1, 4, "19÷7=", "10÷5="),
    @(1, 5, "21÷7=", "19÷8="),

    @(5, 1, "21÷4=", "18÷8="),
    @(5, 2, "66÷5=", "33÷3="),
    @(5, 3, "89÷2=", "54÷9="),
    @(5, 4, "85÷8=", "49÷7="),
    @(5, 5, "45÷3=", "77÷5="),

    @(9, 1, "70÷9=", "86÷9="),
    @(9, 2, "89÷5=", "69÷2="),
    @(9, 3, "52÷2=", "66÷9="),
    @(9, 4, "88÷5=", "74÷4="),
    @(9, 5, "87÷9=", "95÷4="),

    @(13, 1, "36÷5=", "81÷3="),
    @(13, 2, "76÷7=", "10÷4="),
    @(13, 3, "15÷3=", "24÷9="),
    @(13, 4, "89÷5=", "59÷4="),
    @(13, 5, "65÷3=", "36÷4="),

    @(17, 1, "58÷4=", "77÷9="),
    @(17, 2, "58÷4=", "46÷8="),
    @(17, 3, "65÷5=", "30÷7="),
    @(17, 4, "29÷4=", "87÷4="),
    @(17, 5, "51÷4=", "56÷2=")
)

foreach ($edit in $edits) {
    $row = $edit[0]
    $col = $edit[1]
    $old = $edit[2]
    $new = $edit[3]

    $cell = $t.Cell($row, $col)
    $cell.Range.Find.Execute($old, $true, $false, $false, $false, $false,
                              $true, 1, $false, $new, 1)
}
